# "Small changes after first test"
# - Reservoir Age (L) / Reservoir Error (M) columns filled in for all data rows on "Age"
# - Age sheet becomes the active/selected sheet (was Metadata), with a new selection
# - Metadata sheet is no longer the tab-selected sheet

$wb = $excel.ActiveWorkbook

$wsAge = $wb.Worksheets.Item("Age")
$wsMeta = $wb.Worksheets.Item("Metadata")

# Fill in the Reservoir Age (yr) and Reservoir Error (+/- yr) columns (L, M) for rows 2-27
for ($r = 2; $r -le 27; $r++) {
    $wsAge.Cells.Item($r, 12).Value = 853
    $wsAge.Cells.Item($r, 13).Value = 31
}

# Metadata keeps its own stored selection (C5), but loses tab-selected status
# once a different sheet becomes active.
$wsMeta.Activate()
$wsMeta.Range("C5").Select()

# Make "Age" the active sheet (tabSelected) with the new selection at M24
$wsAge.Activate()
$wsAge.Range("M24").Select()
